$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the client code value in row 2
$ws.Range("A2").Value = "24681769"

# Remove row 3 (the second verification record) entirely
$ws.Rows("3:3").Delete()

# Move the active selection to H5, matching the post-edit sheet view
$ws.Range("H5").Select()
